# Insert a new row of data above the current row 186 (shifting the
# existing rows 186-194 down to 187-195) and populate it with the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 186:194 down to 187:195 by inserting a new blank row at 186.
$ws.Rows.Item(186).Insert()

# Fill in the new row 186 with the new record's data.
$ws.Range("A186").Value = 4
$ws.Range("B186").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C186").Value = "Los Lagos"
$ws.Range("D186").Value = [DateTime]::ParseExact("2022-01-07", "yyyy-MM-dd", $null)
$ws.Range("E186").Value = 10
$ws.Range("F186").Value = 100112044
$ws.Range("G186").Value = "Perejil"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 120
$ws.Range("K186").Value = 7000
$ws.Range("L186").Value = 7000
$ws.Range("M186").Value = 7000
$ws.Range("N186").Value = "$/docena de atados (3 kilos)"
$ws.Range("O186").Value = "Región Metropolitana"
$ws.Range("P186").Value = 2333
$ws.Range("Q186").Value = 3
$ws.Range("R186").Value = "Hortaliza"
